$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new numeric-looking text must stay text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "61.667.18"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "3.451.94"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "582.43"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "147.11"
$ws.Range("E6").Value = "  +7.66%  "
$ws.Range("D7").Value = "3.452.55"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "4.039.28"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "27.99"
$ws.Range("E14").Value = "  +10.18%  "
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "3.467.22"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "61.747.58"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +8.92%  "
$ws.Range("D20").Value = "14.39"
$ws.Range("E20").Value = "  +4.14%  "
$ws.Range("D21").Value = "9.57"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "387.54"
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").Value = "73.31"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "5.78"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").Value = "3.605.41"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "7.75"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -9.08%  "
$ws.Range("D33").Value = "8.21"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "24.22"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "3.481.40"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").Value = "7.01"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "166.84"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "0.0787"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("D43").Value = "27.14"
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("D44").Value = "0.808"
$ws.Range("E44").Value = "  +4.43%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "42.52"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "4.51"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "1.73"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "2.573.58"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("E51").Value = "  +2.35%  "
